# naruto stock and image update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stock / price tweaks ---
$ws.Cells.Item(78, 2).Value = 1      # B78: 2 -> 1
$ws.Cells.Item(78, 3).Value = 300    # C78: 270 -> 300

$ws.Cells.Item(81, 2).Value = 2      # B81: 3 -> 2

$ws.Cells.Item(132, 2).Value = 1     # B132: 2 -> 1

$ws.Cells.Item(138, 2).Value = 0     # B138: 1 -> 0

$ws.Cells.Item(142, 3).Value = 220   # C142: 200 -> 220
$ws.Cells.Item(143, 3).Value = 250   # C143: 220 -> 250
$ws.Cells.Item(144, 3).Value = 220   # C144: 200 -> 220
$ws.Cells.Item(147, 3).Value = 220   # C147: 200 -> 220

$ws.Cells.Item(148, 2).Value = 2     # B148: 3 -> 2

$ws.Cells.Item(177, 2).Value = 0     # B177: 1 -> 0
$ws.Cells.Item(178, 2).Value = 0     # B178: 1 -> 0

# --- New product images (assigned in the order they were first introduced
#     so the new shared-string table entries line up 320..333) ---
$ws.Cells.Item(151, 4).Value = "sasuke sus.jpg"
$ws.Cells.Item(152, 4).Value = "madara sus.jpg"
$ws.Cells.Item(153, 4).Value = "kakashi sus.jpg"
$ws.Cells.Item(154, 4).Value = "shisui sus.jpg"
$ws.Cells.Item(155, 4).Value = "itachi sus.jpg"
$ws.Cells.Item(156, 4).Value = "hagoromo sus.jpg"
$ws.Cells.Item(138, 4).Value = "kakashi cape.jpg"
$ws.Cells.Item(140, 4).Value = "madara six path.jpg"
$ws.Cells.Item(141, 4).Value = "madara.jpg"
$ws.Cells.Item(142, 4).Value = "minato.jpg"
$ws.Cells.Item(143, 4).Value = "sasuke.jpg"
$ws.Cells.Item(149, 4).Value = "obito akatsuki.jpg"
$ws.Cells.Item(147, 4).Value = "obito war.jpg"
$ws.Cells.Item(139, 4).Value = "obito six path.jpg"

# --- View / selection state ---
$selected = $ws.Range("D139").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 128
$win.ScrollColumn = 1
